$wb = $excel.ActiveWorkbook

# The same "想去人数" (want-to-go count) updates apply to both the
# "展览" and "全部类型" sheets, which mirror each other's data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F9").Value = 662
    $ws.Range("F19").Value = 3842
    $ws.Range("F24").Value = 710
    $ws.Range("F25").Value = 481
    $ws.Range("F28").Value = 1636
}
